$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Axis Pattern")
Write-Host ("D381 formula before: " + $ws.Cells.Item(381,4).Formula)
$ws.Rows("344:410").Insert()
Write-Host ("D448 formula after: " + $ws.Cells.Item(448,4).Formula)
